$d = $word.ActiveDocument

function Set-ParagraphRuns {
    param($ParaIndex, $RunsXml)

    $para = $d.Paragraphs($ParaIndex)
    $r = $para.Range
    $start = $r.Start
    $end = $r.End

    # Range covering just the run content (exclude the trailing paragraph mark)
    $body = $d.Range($start, $end - 1)
    if ($body.Start -ne $body.End) {
        $body.Delete()
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $ins = $d.Range($start, $start)
    $ins.InsertXML($xml)
}

# --- Paragraph 18: "In Eclipse, import the project directly from hands-on-exercises/Exam Results App/step 2"
$p18Runs = '<w:r><w:t xml:space="preserve">In Eclipse, import the project directly from </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>hands-on-exercises/Exam</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Results</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> A</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>pp/</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>ExamResults</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>-S</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>tep 2</w:t></w:r>'

Set-ParagraphRuns 18 $p18Runs

# --- Paragraph 22: "In the Import Projects dialog, Browse to the folder hands-on-exercises/ExamResultsapp/step 2 and tab out"
$p22Runs = '<w:r><w:t xml:space="preserve">In the </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Import Projects </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">dialog, Browse to the folder </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>hands-on-exercises/Exam</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Results</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> A</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>pp/</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>ExamResults</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>-S</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>tep 2</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t>and tab out</w:t></w:r>'

Set-ParagraphRuns 22 $p22Runs

# --- Paragraph 106: "}" (remove <w:lastRenderedPageBreak/>)
$p106Runs = '<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r>'

Set-ParagraphRuns 106 $p106Runs
